$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the action name for row 11 (Door / Push door -> Force door)
$ws.Range("B11").Value = "Force door"

# Add a violent score for forcing the door to balance scoring
$ws.Range("C11").Value = 2

# Update the active selection on the sheet
$ws.Range("D11").Select()
